$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: straightforward odds updates (no row shift) ---
$ws.Cells.Item(5, 39).Value = 800
$ws.Cells.Item(6, 10).Value = 2.7
$ws.Cells.Item(6, 12).Value = 3.7
$ws.Cells.Item(6, 22).Value = 1.91
$ws.Cells.Item(6, 23).Value = 7.5
$ws.Cells.Item(6, 24).Value = 10.5
$ws.Cells.Item(6, 27).Value = 17.5
$ws.Cells.Item(6, 28).Value = 28
$ws.Cells.Item(6, 33).Value = 9.75
$ws.Cells.Item(6, 37).Value = 29
$ws.Cells.Item(6, 40).Value = 4.1
$ws.Cells.Item(6, 41).Value = 10.75
$ws.Cells.Item(6, 42).Value = 17.5
$ws.Cells.Item(6, 43).Value = 40
$ws.Cells.Item(6, 44).Value = 65
$ws.Cells.Item(6, 46).Value = 2.62
$ws.Cells.Item(6, 50).Value = 17.5
$ws.Cells.Item(6, 53).Value = 110
$ws.Cells.Item(7, 8).Value = 4.45
$ws.Cells.Item(7, 9).Value = 7.4
$ws.Cells.Item(7, 10).Value = 1.82
$ws.Cells.Item(7, 11).Value = 2.4
$ws.Cells.Item(7, 12).Value = 6.6
$ws.Cells.Item(7, 19).Value = 1.28
$ws.Cells.Item(7, 20).Value = 3.46
$ws.Cells.Item(7, 21).Value = 1.75
$ws.Cells.Item(7, 22).Value = 1.85
$ws.Cells.Item(7, 23).Value = 8
$ws.Cells.Item(7, 24).Value = 7.2
$ws.Cells.Item(7, 26).Value = 9.25
$ws.Cells.Item(7, 27).Value = 10.5
$ws.Cells.Item(7, 28).Value = 23
$ws.Cells.Item(7, 30).Value = 9
$ws.Cells.Item(7, 31).Value = 18
$ws.Cells.Item(7, 33).Value = 21
$ws.Cells.Item(7, 34).Value = 55
$ws.Cells.Item(7, 35).Value = 23
$ws.Cells.Item(7, 37).Value = 80
$ws.Cells.Item(7, 38).Value = 65
$ws.Cells.Item(7, 42).Value = 14.5
$ws.Cells.Item(7, 43).Value = 16
$ws.Cells.Item(7, 44).Value = 37
$ws.Cells.Item(7, 47).Value = 7.8
$ws.Cells.Item(7, 49).Value = 8.5
$ws.Cells.Item(7, 50).Value = 40
$ws.Cells.Item(7, 51).Value = 40
$ws.Cells.Item(7, 52).Value = 300
$ws.Cells.Item(7, 53).Value = 300
$ws.Cells.Item(7, 54).Value = 500
$ws.Cells.Item(11, 17).Value = 2
$ws.Cells.Item(11, 18).Value = 1.8
$ws.Cells.Item(21, 17).Value = 2.4
$ws.Cells.Item(21, 18).Value = 1.53

# --- Part 2: insert new row 24 (new fixture) ---
$ws.Rows.Item(24).Insert()

# --- Part 3: populate new row 24 ---
$ws.Cells.Item(24, 1).Value = 'jwH10NVQ'
$ws.Cells.Item(24, 2).Value = "'12/10/2024"
$ws.Cells.Item(24, 3).Value = '23:00'
$ws.Cells.Item(24, 4).Value = 'USA - USL CHAMPIONSHIP'
$ws.Cells.Item(24, 5).Value = 'Oakland Roots'
$ws.Cells.Item(24, 6).Value = 'Phoenix Rising'
$ws.Cells.Item(24, 7).Value = 2.32
$ws.Cells.Item(24, 8).Value = 3.2
$ws.Cells.Item(24, 9).Value = 2.8
$ws.Cells.Item(24, 10).Value = 2.92
$ws.Cells.Item(24, 11).Value = 2.12
$ws.Cells.Item(24, 12).Value = 3.4
$ws.Cells.Item(24, 13).Value = 1.07
$ws.Cells.Item(24, 14).Value = 7
$ws.Cells.Item(24, 15).Value = 1.33
$ws.Cells.Item(24, 16).Value = 3.05
$ws.Cells.Item(24, 17).Value = 2
$ws.Cells.Item(24, 18).Value = 1.75
$ws.Cells.Item(24, 19).Value = 1.4
$ws.Cells.Item(24, 20).Value = 2.75
$ws.Cells.Item(24, 21).Value = 1.78
$ws.Cells.Item(24, 22).Value = 1.93
$ws.Cells.Item(24, 23).Value = 7.7
$ws.Cells.Item(24, 24).Value = 11.25
$ws.Cells.Item(24, 25).Value = 9.25
$ws.Cells.Item(24, 26).Value = 24
$ws.Cells.Item(24, 27).Value = 19.5
$ws.Cells.Item(24, 28).Value = 30
$ws.Cells.Item(24, 29).Value = 7
$ws.Cells.Item(24, 30).Value = 6.3
$ws.Cells.Item(24, 31).Value = 14
$ws.Cells.Item(24, 32).Value = 65
$ws.Cells.Item(24, 33).Value = 8.75
$ws.Cells.Item(24, 34).Value = 14
$ws.Cells.Item(24, 35).Value = 10.25
$ws.Cells.Item(24, 36).Value = 35
$ws.Cells.Item(24, 37).Value = 24
$ws.Cells.Item(24, 38).Value = 35
$ws.Cells.Item(24, 39).Value = 500
$ws.Cells.Item(24, 40).Value = 4.3
$ws.Cells.Item(24, 41).Value = 12.5
$ws.Cells.Item(24, 42).Value = 20
$ws.Cells.Item(24, 43).Value = 50
$ws.Cells.Item(24, 44).Value = 80
$ws.Cells.Item(24, 45).Value = 250
$ws.Cells.Item(24, 46).Value = 2.75
$ws.Cells.Item(24, 47).Value = 7
$ws.Cells.Item(24, 48).Value = 60
$ws.Cells.Item(24, 49).Value = 4.8
$ws.Cells.Item(24, 50).Value = 15.5
$ws.Cells.Item(24, 51).Value = 23
$ws.Cells.Item(24, 52).Value = 70
$ws.Cells.Item(24, 53).Value = 100
$ws.Cells.Item(24, 54).Value = 300
$ws.Cells.Item(24, 55).Value = 51
$ws.Cells.Item(24, 56).Value = 51

# --- Part 4: fix up cells in shifted rows 25 and 28 that changed beyond the pure shift ---
$ws.Cells.Item(25, 8).Value = 3.8
$ws.Cells.Item(25, 9).Value = 4.5
$ws.Cells.Item(25, 10).Value = 2.2
$ws.Cells.Item(25, 12).Value = 4.7
$ws.Cells.Item(25, 13).Value = 1.04
$ws.Cells.Item(25, 14).Value = 8.25
$ws.Cells.Item(25, 16).Value = 3.7
$ws.Cells.Item(25, 20).Value = 2.92
$ws.Cells.Item(25, 24).Value = 8
$ws.Cells.Item(25, 28).Value = 24
$ws.Cells.Item(25, 29).Value = 8.25
$ws.Cells.Item(25, 30).Value = 7.6
$ws.Cells.Item(25, 34).Value = 27
$ws.Cells.Item(25, 36).Value = 75
$ws.Cells.Item(25, 37).Value = 40
$ws.Cells.Item(25, 46).Value = 2.92
$ws.Cells.Item(25, 49).Value = 6.3
$ws.Cells.Item(25, 50).Value = 25
$ws.Cells.Item(25, 54).Value = 350
$ws.Cells.Item(28, 8).Value = 2.87
$ws.Cells.Item(28, 9).Value = 3.7
$ws.Cells.Item(28, 10).Value = 2.75
$ws.Cells.Item(28, 12).Value = 4.4
$ws.Cells.Item(28, 21).Value = 2
$ws.Cells.Item(28, 22).Value = 1.65
$ws.Cells.Item(28, 23).Value = 5.8
$ws.Cells.Item(28, 24).Value = 9.25
$ws.Cells.Item(28, 25).Value = 9
$ws.Cells.Item(28, 27).Value = 20
$ws.Cells.Item(28, 28).Value = 37
$ws.Cells.Item(28, 29).Value = 6.3
$ws.Cells.Item(28, 33).Value = 7.9
$ws.Cells.Item(28, 34).Value = 18
$ws.Cells.Item(28, 35).Value = 13.5
$ws.Cells.Item(28, 38).Value = 65
$ws.Cells.Item(28, 41).Value = 11.25
$ws.Cells.Item(28, 42).Value = 22
$ws.Cells.Item(28, 43).Value = 45
$ws.Cells.Item(28, 44).Value = 90
$ws.Cells.Item(28, 47).Value = 7.6
$ws.Cells.Item(28, 48).Value = 90
$ws.Cells.Item(28, 49).Value = 5.4
$ws.Cells.Item(28, 50).Value = 23
$ws.Cells.Item(28, 51).Value = 35
$ws.Cells.Item(28, 52).Value = 150
